$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "40.454.00"
$ws.Range("E2").Value = "  -3.12%  "

$ws.Range("D3").Value = "2.366.36"
$ws.Range("E3").Value = "  -4.56%  "

$ws.Range("E4").Value = "  +0.15%  "

Set-TextValue $ws.Range("D5") "311.18"
$ws.Range("E5").Value = "  -2.46%  "

Set-TextValue $ws.Range("D6") "86.74"
$ws.Range("E6").Value = "  -7.51%  "

Set-TextValue $ws.Range("D7") "0.528"
$ws.Range("E7").Value = "  -4.71%  "

$ws.Range("E8").Value = "  +0.19%  "

Set-TextValue $ws.Range("D9") "0.487"
$ws.Range("E9").Value = "  -6.23%  "

Set-TextValue $ws.Range("D10") "0.0818"
$ws.Range("E10").Value = "  -5.34%  "

Set-TextValue $ws.Range("D11") "30.78"
$ws.Range("E11").Value = "  -7.68%  "

Set-TextValue $ws.Range("D12") "0.108"
$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("D13").Value = "2.748.84"
$ws.Range("E13").Value = "  -3.89%  "

Set-TextValue $ws.Range("D14") "6.58"
$ws.Range("E14").Value = "  -5.01%  "

Set-TextValue $ws.Range("D15") "15.02"
$ws.Range("E15").Value = "  -5.06%  "

$ws.Range("D16").Value = "2.394.40"
$ws.Range("E16").Value = "  -3.52%  "

Set-TextValue $ws.Range("D17") "0.752"
$ws.Range("E17").Value = "  -5.09%  "

$ws.Range("D18").Value = "40.487.18"
$ws.Range("E18").Value = "  -2.93%  "

$ws.Range("D19").Value = "0.0₃0906"
$ws.Range("E19").Value = "  -4.93%  "

Set-TextValue $ws.Range("D20") "6.11"
$ws.Range("E20").Value = "  -5.74%  "

Set-TextValue $ws.Range("D21") "68.72"
$ws.Range("E21").Value = "  -3.63%  "

Set-TextValue $ws.Range("D22") "10.65"
$ws.Range("E22").Value = "  -6.51%  "

Set-TextValue $ws.Range("D23") "234.07"
$ws.Range("E23").Value = "  -2.56%  "

Set-TextValue $ws.Range("D24") "2.62"
$ws.Range("E24").Value = "  -4.93%  "

$ws.Range("E25").Value = "  -0.01%  "

Set-TextValue $ws.Range("D26") "1.81"
$ws.Range("E26").Value = "  -7.03%  "

Set-TextValue $ws.Range("D27") "23.46"
$ws.Range("E27").Value = "  -5.45%  "

Set-TextValue $ws.Range("D28") "2.20"
$ws.Range("E28").Value = "  -2.73%  "

Set-TextValue $ws.Range("D29") "9.32"
$ws.Range("E29").Value = "  -5.35%  "

Set-TextValue $ws.Range("D30") "33.58"
$ws.Range("E30").Value = "  -7.35%  "

Set-TextValue $ws.Range("D31") "154.83"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("E32").Value = "  +0.11%  "

Set-TextValue $ws.Range("D33") "5.19"
$ws.Range("E33").Value = "  -6.45%  "

Set-TextValue $ws.Range("D34") "0.0725"
$ws.Range("E34").Value = "  -6.02%  "

$ws.Range("E35").Value = "  -6.31%  "

Set-TextValue $ws.Range("D36") "0.113"
$ws.Range("E36").Value = "  -2.16%  "

Set-TextValue $ws.Range("D37") "2.80"
$ws.Range("E37").Value = "  -4.92%  "

Set-TextValue $ws.Range("D38") "15.71"
$ws.Range("E38").Value = "  -9.78%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D39") "0.0982"
$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D40") "1.71"
$ws.Range("E40").Value = "  -9.38%  "

Set-TextValue $ws.Range("D41") "3.79"
$ws.Range("E41").Value = "  -6.51%  "

$ws.Range("E42").Value = "  -7.46%  "

$ws.Range("D43").Value = "1.953.38"
$ws.Range("E43").Value = "  -2.16%  "

Set-TextValue $ws.Range("D44") "0.0268"
$ws.Range("E44").Value = "  -6.34%  "

Set-TextValue $ws.Range("D45") "17.62"
$ws.Range("E45").Value = "  -9.63%  "

Set-TextValue $ws.Range("D46") "2.78"
$ws.Range("E46").Value = "  -7.59%  "

Set-TextValue $ws.Range("D47") "9.24"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").Value = "2.621.21"
$ws.Range("E48").Value = "  -3.47%  "

Set-TextValue $ws.Range("D49") "72.40"
$ws.Range("E49").Value = "  -2.86%  "

Set-TextValue $ws.Range("D50") "92.88"
$ws.Range("E50").Value = "  -4.89%  "

Set-TextValue $ws.Range("D51") "49.73"
$ws.Range("E51").Value = "  -5.62%  "
